$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: O11 changes from 572880.59 to 624605.1
$ws.Range("O11").Value = 624605.1

# Row 14: O14 changes from 1832.18 to 2073.75
$ws.Range("O14").Value = 2073.75

# Row 16: N16 and O16 were empty, now both get 701.99
$ws.Range("N16").Value = 701.99
$ws.Range("O16").Value = 701.99

# Row 28: N28 changes from 57376.24 to 63376.24
$ws.Range("N28").Value = 63376.24
